$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 changes from the text "R40" to the text "1". A plain
# Value assignment of "1" would be auto-detected as a number, but the
# target keeps the cell as a (shared) text string, so we stage a
# formula in a scratch cell that evaluates to the text "1", copy it,
# and paste-special "values only" into B11. That way the destination
# keeps its existing style/formatting (borders, fill, etc.) untouched
# and no number-auto-detection kicks in.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=TEXT(1,""0"")"
$scratch.Copy()

$target = $ws.Range("B11")
$target.PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
$excel.CutCopyMode = $false
